# Refresh the "cryptos" price/volume snapshot (GitHub Actions scheduled update).
# Price (column D) and 1h volume-change (column E) values are refreshed for every
# coin row, and the Hedera / InternetComputer(DFINITY) rows (37-38) swap ranking
# position while keeping their original rank index (column A) untouched.
# D-column values are prefixed with a leading apostrophe so Excel stores them as
# text (matching the source data, which includes values like "1.000" / "30.358.18"
# that must not be reinterpreted as numbers/dates).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.360.93"
$ws.Range("E2").Value = "  +1.34%  "
$ws.Range("D3").Value = "'2.011.47"
$ws.Range("E3").Value = "  +4.87%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'325.23"
$ws.Range("E5").Value = "  +1.55%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").Value = "'0.5121"
$ws.Range("E7").Value = "  +1.33%  "
$ws.Range("D8").Value = "'0.4263"
$ws.Range("E8").Value = "  +5.43%  "
$ws.Range("D9").Value = "'0.08716"
$ws.Range("E9").Value = "  +4.53%  "
$ws.Range("D10").Value = "'43.31"
$ws.Range("E10").Value = "  +3.21%  "
$ws.Range("E11").Value = "  +2.92%  "
$ws.Range("D12").Value = "'24.61"
$ws.Range("E12").Value = "  +2.86%  "
$ws.Range("D13").Value = "'2.008.10"
$ws.Range("E13").Value = "  +4.48%  "
$ws.Range("D14").Value = "'6.590"
$ws.Range("E14").Value = "  +2.92%  "
$ws.Range("D15").Value = "'7.458"
$ws.Range("E15").Value = "  +2.98%  "
$ws.Range("D16").Value = "'1.001"
$ws.Range("E16").Value = "  -0.31%  "
$ws.Range("D17").Value = "'94.22"
$ws.Range("E17").Value = "  +2.11%  "
$ws.Range("D18").Value = "'0.00001113"
$ws.Range("E18").Value = "  +1.31%  "
$ws.Range("D19").Value = "'0.06521"
$ws.Range("E19").Value = "  +0.13%  "
$ws.Range("E20").Value = "  +3.48%  "
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").Value = "'6.206"
$ws.Range("E22").Value = "  +4.31%  "
$ws.Range("D23").Value = "'30.421.60"
$ws.Range("E23").Value = "  +1.42%  "
$ws.Range("D24").Value = "'11.83"
$ws.Range("E24").Value = "  +4.58%  "
$ws.Range("D25").Value = "'2.265"
$ws.Range("E25").Value = "  +3.07%  "
$ws.Range("D26").Value = "'2.252.37"
$ws.Range("E26").Value = "  +5.16%  "
$ws.Range("D27").Value = "'22.43"
$ws.Range("E27").Value = "  +1.42%  "
$ws.Range("D28").Value = "'162.42"
$ws.Range("E28").Value = "  -0.13%  "
$ws.Range("D29").Value = "'2.431"
$ws.Range("E29").Value = "  +4.88%  "
$ws.Range("D30").Value = "'130.96"
$ws.Range("E30").Value = "  +1.25%  "
$ws.Range("D31").Value = "'1.142"
$ws.Range("E31").Value = "  +0.94%  "
$ws.Range("E32").Value = "  +1.56%  "
$ws.Range("D33").Value = "'6.091"
$ws.Range("E33").Value = "  +2.18%  "
$ws.Range("D34").Value = "'3.827"
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("D35").Value = "'1.374"
$ws.Range("E35").Value = "  +14.82%  "
$ws.Range("D36").Value = "'0.02525"
$ws.Range("E36").Value = "  +2.90%  "
$ws.Range("B37").Value = "InternetComputer(DFINITY)"
$ws.Range("C37").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D37").Value = "'5.458"
$ws.Range("E37").Value = "  +1.00%  "
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "'0.06669"
$ws.Range("E38").Value = "  +3.92%  "
$ws.Range("D39").Value = "'12.35"
$ws.Range("E39").Value = "  +8.71%  "
$ws.Range("D40").Value = "'9.140"
$ws.Range("E40").Value = "  +4.76%  "
$ws.Range("D41").Value = "'0.2194"
$ws.Range("E41").Value = "  +2.01%  "
$ws.Range("E42").Value = "  +2.93%  "
$ws.Range("D43").Value = "'1.239"
$ws.Range("E43").Value = "  +2.45%  "
$ws.Range("D44").Value = "'1.000"
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("D45").Value = "'13.65"
$ws.Range("E45").Value = "  +2.30%  "
$ws.Range("D46").Value = "'0.6162"
$ws.Range("E46").Value = "  +2.05%  "
$ws.Range("D47").Value = "'2.195"
$ws.Range("E47").Value = "  -1.14%  "
$ws.Range("D48").Value = "'3.670"
$ws.Range("E48").Value = "  +0.88%  "
$ws.Range("D49").Value = "'1.261"
$ws.Range("E49").Value = "  +4.25%  "
$ws.Range("D50").Value = "'124.60"
$ws.Range("E50").Value = "  +1.87%  "
$ws.Range("D51").Value = "'80.69"
$ws.Range("E51").Value = "  +2.15%  "
